$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added at the top of the data table (row 27),
# pushing every existing data row down by one (old row 27 -> 28, ... old row
# 130 -> 131). Insert a fresh row at 27 to shift everything down, then
# populate it with the new record's values.
$ws.Rows.Item(27).Insert()

$ws.Range("A27").Value = 11
$ws.Range("B27").Value = "Vega Monumental Concepción"
$ws.Range("C27").Value = "Bíobío"
$ws.Range("D27").Value = 45133
$ws.Range("E27").Value = 8
$ws.Range("F27").Value = 100112012
$ws.Range("G27").Value = "Espinaca"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = 7000
$ws.Range("L27").Value = 7500
$ws.Range("M27").Value = 7200
$ws.Range("N27").Value = "$/cuna 10 kilos"
$ws.Range("O27").Value = "Región Metropolitana"
$ws.Range("P27").Value = 720
$ws.Range("Q27").Value = 10
$ws.Range("R27").Value = "Hortaliza"
